# "allowed camp committee members to generate reports"
#
# The staff list keeps each person's report-access code in column E
# ("Password"). Madhukumar (row 2, SCSE faculty / camp committee) had
# code "2"; bump it to "12" so the camp committee can generate reports.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("E2")

# Remember the cell's current style so the text-entry below (which would
# otherwise be auto-detected as the number 12) doesn't pick up a new
# number format / quote-prefix style - we just want a plain text value,
# same as the "2" it replaces.
$origStyle = $cell.Style

# Leading apostrophe forces this to be stored as text ("12"), not the
# numeric value 12, matching how the rest of column E stores its codes.
$cell.Value = "'12"

$cell.Style = $origStyle
